# Apply updated crypto market data (prices & 1h volume change) per the
# Mon Jul  8 07:08:25 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.574.61"
$ws.Range("E2").Value = "  -3.64%  "
$ws.Range("D3").Value = "2.907.75"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.96"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.06"
$ws.Range("E6").Value = "  -5.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.420"
$ws.Range("E8").Value = "  -4.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.13"
$ws.Range("E9").Value = "  -5.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  -6.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.350"
$ws.Range("E11").Value = "  -5.05%  "
$ws.Range("D12").Value = "3.420.19"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.62"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000159"
$ws.Range("E15").Value = "  -4.78%  "
$ws.Range("D16").Value = "55.525.60"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.95"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("D18").Value = "2.912.34"
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.68"
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.484"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.83"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "3.040.46"
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -6.01%  "
$ws.Range("D28").Value = "0.0₃0831"
$ws.Range("E28").Value = "  -11.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.27"
$ws.Range("E29").Value = "  -8.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.77"
$ws.Range("E30").Value = "  -9.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.76"
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.13"
$ws.Range("E33").Value = "  -6.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.40"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.35"
$ws.Range("E35").Value = "  -8.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.57"
$ws.Range("E36").Value = "  -5.63%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -8.26%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.55"
$ws.Range("E38").Value = "  -5.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0641"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.33"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.635"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.97"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.104.52"
$ws.Range("E45").Value = "  -9.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.32"
$ws.Range("E46").Value = "  -6.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.916"
$ws.Range("E47").Value = "  -8.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0233"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.54"
$ws.Range("E49").Value = "  -5.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0835"
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.65"
$ws.Range("E51").Value = "  -11.50%  "
